$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-01-11 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-12 Friday", 2) | Out-Null

# Update table cell values (row-major order), using Table.Cell(row, col) to
# avoid ambiguity from duplicate cell text (e.g. "30-0=30" appears twice).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "73-73=0"
$t.Cell(1, 2).Range.Text = "90-29=61"
$t.Cell(1, 3).Range.Text = "3+6=9"
$t.Cell(1, 4).Range.Text = "17+56=73"
$t.Cell(1, 5).Range.Text = "45+38=83"
$t.Cell(2, 1).Range.Text = "5+79=84"
$t.Cell(2, 2).Range.Text = "4+2=6"
$t.Cell(2, 3).Range.Text = "78-0=78"
$t.Cell(2, 4).Range.Text = "63+5=68"
$t.Cell(2, 5).Range.Text = "49+41=90"
$t.Cell(3, 1).Range.Text = "44+34=78"
$t.Cell(3, 2).Range.Text = "55-5=50"
$t.Cell(3, 3).Range.Text = "77+15=92"
$t.Cell(3, 4).Range.Text = "2+62=64"
$t.Cell(3, 5).Range.Text = "89-49=40"
$t.Cell(4, 1).Range.Text = "65-51=14"
$t.Cell(4, 2).Range.Text = "79-46=33"
$t.Cell(4, 3).Range.Text = "55-47=8"
$t.Cell(4, 4).Range.Text = "54+4=58"
$t.Cell(4, 5).Range.Text = "63-58=5"
$t.Cell(5, 1).Range.Text = "37+42=79"
$t.Cell(5, 2).Range.Text = "2+39=41"
$t.Cell(5, 3).Range.Text = "16-12=4"
$t.Cell(5, 4).Range.Text = "70+10=80"
$t.Cell(5, 5).Range.Text = "22+63=85"
$t.Cell(6, 1).Range.Text = "93-18=75"
$t.Cell(6, 2).Range.Text = "10+20=30"
$t.Cell(6, 3).Range.Text = "39-34=5"
$t.Cell(6, 4).Range.Text = "54+40=94"
$t.Cell(6, 5).Range.Text = "92-45=47"
$t.Cell(7, 1).Range.Text = "54-28=26"
$t.Cell(7, 2).Range.Text = "33+22=55"
$t.Cell(7, 3).Range.Text = "76+15=91"
$t.Cell(7, 4).Range.Text = "96-90=6"
$t.Cell(7, 5).Range.Text = "0+85=85"
$t.Cell(8, 1).Range.Text = "71+27=98"
$t.Cell(8, 2).Range.Text = "67+7=74"
$t.Cell(8, 3).Range.Text = "95+4=99"
$t.Cell(8, 4).Range.Text = "76-38=38"
$t.Cell(8, 5).Range.Text = "40+29=69"
$t.Cell(9, 1).Range.Text = "45+24=69"
$t.Cell(9, 2).Range.Text = "28+35=63"
$t.Cell(9, 3).Range.Text = "69-21=48"
$t.Cell(9, 4).Range.Text = "50+36=86"
$t.Cell(9, 5).Range.Text = "26-26=0"
$t.Cell(10, 1).Range.Text = "9+82=91"
$t.Cell(10, 2).Range.Text = "14+3=17"
$t.Cell(10, 3).Range.Text = "68+22=90"
$t.Cell(10, 4).Range.Text = "34-6=28"
$t.Cell(10, 5).Range.Text = "35+59=94"
$t.Cell(11, 1).Range.Text = "18+17=35"
$t.Cell(11, 2).Range.Text = "84-41=43"
$t.Cell(11, 3).Range.Text = "88+1=89"
$t.Cell(11, 4).Range.Text = "97-87=10"
$t.Cell(11, 5).Range.Text = "69+8=77"
$t.Cell(12, 1).Range.Text = "2+25=27"
$t.Cell(12, 2).Range.Text = "11+77=88"
$t.Cell(12, 3).Range.Text = "95-61=34"
$t.Cell(12, 4).Range.Text = "78+9=87"
$t.Cell(12, 5).Range.Text = "76-49=27"
$t.Cell(13, 1).Range.Text = "93-88=5"
$t.Cell(13, 2).Range.Text = "65-26=39"
$t.Cell(13, 3).Range.Text = "86-46=40"
$t.Cell(13, 4).Range.Text = "4+61=65"
$t.Cell(13, 5).Range.Text = "3+85=88"
$t.Cell(14, 1).Range.Text = "85-1=84"
$t.Cell(14, 2).Range.Text = "55-43=12"
$t.Cell(14, 3).Range.Text = "75-50=25"
$t.Cell(14, 4).Range.Text = "86+12=98"
$t.Cell(14, 5).Range.Text = "20-11=9"
$t.Cell(15, 1).Range.Text = "33+7=40"
$t.Cell(15, 2).Range.Text = "78-71=7"
$t.Cell(15, 3).Range.Text = "31+5=36"
$t.Cell(15, 4).Range.Text = "37+56=93"
$t.Cell(15, 5).Range.Text = "37+9=46"
$t.Cell(16, 1).Range.Text = "26+9=35"
$t.Cell(16, 2).Range.Text = "29+50=79"
$t.Cell(16, 3).Range.Text = "72-11=61"
$t.Cell(16, 4).Range.Text = "95-6=89"
$t.Cell(16, 5).Range.Text = "35+12=47"
$t.Cell(17, 1).Range.Text = "67-60=7"
$t.Cell(17, 2).Range.Text = "59+18=77"
$t.Cell(17, 3).Range.Text = "86-58=28"
$t.Cell(17, 4).Range.Text = "57+32=89"
$t.Cell(17, 5).Range.Text = "47+4=51"
$t.Cell(18, 1).Range.Text = "49-33=16"
$t.Cell(18, 2).Range.Text = "64-42=22"
$t.Cell(18, 3).Range.Text = "50+49=99"
$t.Cell(18, 4).Range.Text = "95-26=69"
$t.Cell(18, 5).Range.Text = "8+24=32"
$t.Cell(19, 1).Range.Text = "94-48=46"
$t.Cell(19, 2).Range.Text = "32+67=99"
$t.Cell(19, 3).Range.Text = "48+41=89"
$t.Cell(19, 4).Range.Text = "57-26=31"
$t.Cell(19, 5).Range.Text = "19+39=58"
$t.Cell(20, 1).Range.Text = "63+2=65"
$t.Cell(20, 2).Range.Text = "98-67=31"
$t.Cell(20, 3).Range.Text = "99-45=54"
$t.Cell(20, 4).Range.Text = "49+20=69"
$t.Cell(20, 5).Range.Text = "18+1=19"
